# Share CIT Nationwide Total history
# Insert per-state regression result blocks (State / Point estimate /
# 95% CI / t-statistic / p-value) right after the "Nationwide Corporate
# income ... doesn't have to sum to 1." paragraph, before the existing
# blank paragraph that follows it.

$d = $word.ActiveDocument

$anchorText = "Nationwide Corporate income for a company doesn" + [char]0x2019 + "t have to sum to 1."

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq $anchorText) {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Anchor paragraph not found"
}

# Add a blank paragraph right after the anchor (mirrors the existing
# blank paragraph pattern already used elsewhere in the doc).
$anchor.Range.InsertParagraphAfter()

# Re-resolve the freshly inserted blank paragraph, then drop all of the
# new "State: ..." blocks right after it in a single InsertAfter call
# (one `r per paragraph break).
$blank = $anchor.Next()
$insertionPoint = $blank.Range
$insertionPoint.Collapse(0)

$lines = @(
    'State: Iowa',
    'Point estimate: 0.00',
    '95% CI (-0.02, 0.02)',
    't-statistic: 0.06',
    'p-value: 0.9551',
    'State: Nebraska',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.09',
    'p-value: 0.9304',
    'State: Michigan',
    'Point estimate: 0.01',
    '95% CI (0.00, 0.01)',
    't-statistic: 2.40',
    'p-value: 0.0166',
    'State: Illinois',
    'Point estimate: 0.01',
    '95% CI (0.01, 0.02)',
    't-statistic: 4.53',
    'p-value: 0.0000',
    'State: Oregon',
    'Point estimate: 0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: 0.09',
    'p-value: 0.9302',
    'State: Georgia',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.00)',
    't-statistic: -1.07',
    'p-value: 0.2856',
    'State: Wisconsin',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.00)',
    't-statistic: -0.26',
    'p-value: 0.7920',
    'State: Arizona',
    'Point estimate: 0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: 0.04',
    'p-value: 0.9689',
    'State: Indiana',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.00)',
    't-statistic: -0.98',
    'p-value: 0.3251',
    'State: Maine',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.26',
    'p-value: 0.7944',
    'State: Minnesota',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.57',
    'p-value: 0.5683',
    'State: Pennsylvania',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.73',
    'p-value: 0.4683',
    'State: South Carolina',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.44',
    'p-value: 0.6599',
    'State: Colorado',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.36',
    'p-value: 0.7197',
    'State: California',
    'Point estimate: -0.04',
    '95% CI (-0.04, -0.03)',
    't-statistic: -10.98',
    'p-value: 0.0000',
    'State: Utah',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.05',
    'p-value: 0.9599',
    'State: New Jersey',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.13',
    'p-value: 0.8984',
    'State: New York',
    'Point estimate: -0.01',
    '95% CI (-0.02, -0.00)',
    't-statistic: -2.93',
    'p-value: 0.0034',
    'State: Rhode Island',
    'Point estimate: 0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: 0.06',
    'p-value: 0.9533',
    'State: Connecticut',
    'Point estimate: 0.00',
    '95% CI (-0.00, 0.01)',
    't-statistic: 1.25',
    'p-value: 0.2110',
    'State: Louisiana',
    'Point estimate: 0.00',
    '95% CI (-0.00, 0.01)',
    't-statistic: 0.39',
    'p-value: 0.6947',
    'State: North Carolina',
    'Point estimate: -0.01',
    '95% CI (-0.01, -0.00)',
    't-statistic: -3.66',
    'p-value: 0.0003',
    'State: North Dakota',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.00)',
    't-statistic: -0.83',
    'p-value: 0.4057',
    'State: Delaware',
    'Point estimate: -0.00',
    '95% CI (-0.00, 0.00)',
    't-statistic: -1.12',
    'p-value: 0.2634',
    'State: Kentucky',
    'Point estimate: 0.00',
    '95% CI (-0.00, 0.01)',
    't-statistic: 0.72',
    'p-value: 0.4735',
    'State: Maryland',
    'Point estimate: 0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: 0.05',
    'p-value: 0.9594',
    'State: Missouri',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.00)',
    't-statistic: -0.37',
    'p-value: 0.7082',
    'State: Alabama',
    'Point estimate: 0.00',
    '95% CI (-0.00, 0.01)',
    't-statistic: 1.02',
    'p-value: 0.3064',
    'State: Arkansas',
    'Point estimate: -0.00',
    '95% CI (-0.01, 0.01)',
    't-statistic: -0.11',
    'p-value: 0.9151'
)

$text = [string]::Join("`r", $lines)
$insertionPoint.InsertAfter($text)

Write-Output "Inserted $($lines.Count) paragraphs."
